$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $c = $sheet.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "28.328.38"
Set-TextValue $ws "E2" "  -5.18%  "
Set-TextValue $ws "D3" "1.841.45"
Set-TextValue $ws "E3" "  -4.84%  "
Set-TextValue $ws "D4" "1.001"
Set-TextValue $ws "E4" "  -0.63%  "
Set-TextValue $ws "D5" "330.14"
Set-TextValue $ws "E5" "  -1.59%  "
Set-TextValue $ws "E6" "  -0.47%  "
Set-TextValue $ws "D7" "0.4603"
Set-TextValue $ws "E7" "  -4.68%  "
Set-TextValue $ws "D8" "0.3866"
Set-TextValue $ws "E8" "  -6.01%  "
Set-TextValue $ws "D9" "46.11"
Set-TextValue $ws "E9" "  -2.67%  "
Set-TextValue $ws "E10" "  -3.55%  "
Set-TextValue $ws "D11" "0.9654"
Set-TextValue $ws "E11" "  -4.79%  "
Set-TextValue $ws "D12" "21.98"
Set-TextValue $ws "E12" "  -7.04%  "
Set-TextValue $ws "D13" "1.809.30"
Set-TextValue $ws "E13" "  -8.09%  "
Set-TextValue $ws "E14" "  -5.91%  "
Set-TextValue $ws "D15" "6.931"
Set-TextValue $ws "E15" "  -4.88%  "
Set-TextValue $ws "D16" "0.06876"
Set-TextValue $ws "E16" "  +0.28%  "
Set-TextValue $ws "D17" "1.001"
Set-TextValue $ws "E17" "  -0.68%  "
Set-TextValue $ws "D18" "86.97"
Set-TextValue $ws "E18" "  -4.52%  "
Set-TextValue $ws "D19" "0.000009951"
Set-TextValue $ws "E19" "  -3.87%  "
Set-TextValue $ws "D20" "16.96"
Set-TextValue $ws "E20" "  -4.69%  "
Set-TextValue $ws "E21" "  -0.35%  "
Set-TextValue $ws "D22" "28.358.75"
Set-TextValue $ws "E22" "  -5.09%  "
Set-TextValue $ws "D23" "5.345"
Set-TextValue $ws "E23" "  -5.03%  "
Set-TextValue $ws "D24" "10.98"
Set-TextValue $ws "E24" "  -7.40%  "
Set-TextValue $ws "D25" "2.136"
Set-TextValue $ws "E25" "  -1.99%  "
Set-TextValue $ws "D26" "2.069.61"
Set-TextValue $ws "E26" "  -6.05%  "
Set-TextValue $ws "D27" "153.70"
Set-TextValue $ws "E27" "  -1.89%  "
Set-TextValue $ws "D28" "19.21"
Set-TextValue $ws "E28" "  -4.16%  "
Set-TextValue $ws "E29" "  -13.63%  "
Set-TextValue $ws "D30" "1.993"
Set-TextValue $ws "E30" "  -4.84%  "
Set-TextValue $ws "D31" "117.16"
Set-TextValue $ws "E31" "  -3.40%  "
Set-TextValue $ws "D32" "0.9438"
Set-TextValue $ws "E32" "  -6.24%  "
Set-TextValue $ws "E33" "  -3.25%  "
Set-TextValue $ws "E34" "  -5.12%  "
Set-TextValue $ws "E35" "  -2.88%  "
Set-TextValue $ws "E36" "  -6.19%  "
Set-TextValue $ws "D37" "0.06019"
Set-TextValue $ws "E37" "  -8.30%  "
Set-TextValue $ws "E38" "  -5.56%  "
Set-TextValue $ws "D39" "1.148"
Set-TextValue $ws "E39" "  -4.70%  "
Set-TextValue $ws "D40" "1.000"
Set-TextValue $ws "E40" "  -0.46%  "
Set-TextValue $ws "D41" "7.649"
Set-TextValue $ws "E41" "  -3.96%  "
Set-TextValue $ws "D42" "0.5636"
Set-TextValue $ws "E42" "  -5.49%  "
Set-TextValue $ws "D43" "10.04"
Set-TextValue $ws "E43" "  -6.38%  "
Set-TextValue $ws "E44" "  -3.37%  "
Set-TextValue $ws "B45" "WEMIXToken"
Set-TextValue $ws "C45" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D45" "1.213"
Set-TextValue $ws "E45" "  -4.40%  "
Set-TextValue $ws "B46" "RenderToken"
Set-TextValue $ws "C46" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D46" "2.273"
Set-TextValue $ws "E46" "  -8.58%  "
Set-TextValue $ws "E47" "  -4.74%  "
Set-TextValue $ws "E48" "  -5.82%  "
Set-TextValue $ws "D49" "0.07048"
Set-TextValue $ws "E49" "  -5.76%  "
Set-TextValue $ws "E50" "  -7.45%  "
Set-TextValue $ws "D51" "112.83"
Set-TextValue $ws "E51" "  -3.72%  "
